$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the "Help Menu" label to lowercase "help menu"
$ws.Range("D5").Value = "help menu"

# Normalize B4 cell style to match the non-duplicate style used elsewhere (no explicit fill)
$ws.Range("B4").Interior.Pattern = -4142

# Update the active selection to D5, matching the cell that was edited
$ws.Range("D5").Select()
